$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The MD (Men's Doubles) roster section (rows 11-22) gains an extra pair,
# shifting every following row (the WD and MS sections) down by one row.
# A duplicate trailing row in the old MS section (rows 44 & 45, both
# "Kai Wagner") collapses back down to a single row, keeping the sheet's
# total row count the same.
# ---------------------------------------------------------------------------

# 1) Insert a new row at 23. This shifts old rows 23-45 down to 24-46,
#    and correctly carries merged cells (B24:C24 -> B25:C25, B33:C33 -> B34:C34, etc.)
$ws.Rows("23").Insert()

# 2) Copy the formatting of row 22 (the last MD row) onto the freshly
#    inserted row 23 so it matches the rest of the MD table, then fix the
#    row height which PasteSpecial(Formats) does not carry over.
$ws.Range("A22:L22").Copy()
$ws.Range("A23:L23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows("23").RowHeight = 15.75

# 3) New row 23 becomes MD pair #13: Victor Shi / Curtis Luu (this was
#    previously row 22's content).
$ws.Range("A23").Value2 = 13.0
$ws.Range("B23").Value2 = "Victor Shi"
$ws.Range("C23").Value2 = "Curtis Luu"

# 4) Row 22 (MD pair #12) now becomes Shakthi Guruswami / Rithwik Vaidun.
$ws.Range("B22").Value2 = "Shakthi Guruswami"
$ws.Range("C22").Value2 = "Rithwik Vaidun"

# 5) Remove the duplicate "Kai Wagner" row that the shift produced at the
#    bottom of the MS section (old row 45 is now row 46).
$ws.Rows("46").Delete()
